$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The two test-case rows now record a "PASSED" result in column F
$ws.Range("F2").Value = "PASSED"
$ws.Range("F3").Value = "PASSED"
